$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 20, pushing existing rows 20-27
# down to 21-28 (matches the diff: a new weekly record is prepended to the
# Primera/Start Ruby block, old rows shift down one position).
$ws.Rows.Item(20).EntireRow.Insert()

# Populate the newly inserted row 20 with the new weekly data point.
$ws.Cells.Item(20, 1).Value = 6
$ws.Cells.Item(20, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(20, 3).Value = "Metropolitana"
$ws.Cells.Item(20, 4).Value = 44609
$ws.Cells.Item(20, 5).Value = 13
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100102
$ws.Cells.Item(20, 8).Value = "Cítricos"
$ws.Cells.Item(20, 9).Value = 100102006
$ws.Cells.Item(20, 10).Value = "Pomelo"
$ws.Cells.Item(20, 11).Value = "Start Ruby"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 24
$ws.Cells.Item(20, 14).Value = 190000
$ws.Cells.Item(20, 15).Value = 190000
$ws.Cells.Item(20, 16).Value = 190000
$ws.Cells.Item(20, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(20, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 19).Value = 543
$ws.Cells.Item(20, 20).Value = 350
